$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Kevin Huerter", "SG,SF", "Sacramento Kings"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Terry Rozier", "PG", "Miami Heat"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Shai Gilgeous-Alexander", "PG,SG", "Oklahoma City Thunder"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Keyonte George", "PG,SG", "Utah Jazz")
)

$row = 2
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $row = $row + 1
}
